# "Generate Report for Handoff" — the localization CI job re-ran and picked
# up a new handoff: status flips from "In Translation" to "Ready for
# handoff" and the handoff timestamps advance a few seconds. The Status
# columns widen to fit the longer "Ready for handoff" text.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Status text: "In Translation" -> "Ready for handoff" ---
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$zhcn.Range("C2").Value     = "Ready for handoff"
$dede.Range("C2").Value     = "Ready for handoff"

# --- Handoff timestamps move forward (new xliff generation run) ---
$zhcn.Range("H2").Value     = "2016-08-28 12:59:25"
$dede.Range("H2").Value     = "2016-08-28 12:59:29"
$overview.Range("G2").Value = "2016-08-28 12:59:29"

# --- Widen the Status columns to fit "Ready for handoff" ---
# A ColumnWidth of 16.333333333333332 "characters" is the input that the
# pixel-grid column-width model (MDW=7) rounds to the stored width closest
# to the target 17.2159881591797 used by the report generator.
$newStatusWidth = 16.333333333333332
$overview.Columns.Item(5).ColumnWidth = $newStatusWidth   # zh-cn status column
$overview.Columns.Item(6).ColumnWidth = $newStatusWidth   # de-de status column
$zhcn.Columns.Item(3).ColumnWidth     = $newStatusWidth   # Status column
$dede.Columns.Item(3).ColumnWidth     = $newStatusWidth   # Status column
